$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ---
$ws.Name = "searchitems"

# --- Populate the search-items list ---
$ws.Range("A1").Value = "SearchItems"
$ws.Range("A2").Value = "T-shirts"
$ws.Range("A3").Value = "shoes"
$ws.Range("A4").Value = "kitchen"
$ws.Range("A5").Value = "television"

# --- Header (A1) formatting: bold font + yellow fill ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535

# --- Thin border around every populated cell (A1:A5) ---
$ws.Range("A1:A5").Borders.LineStyle = 1

# --- Column width (emits as width=19 once Excel's default cell-padding is backed out) ---
$ws.Columns("A").ColumnWidth = 18.16666666666667

# --- Selection shown when the sheet is reopened ---
$ws.Range("C5").Select() | Out-Null

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
